# Flow Focusing Q calculator v1.0.0
# Appends a new data row (row 13) to Sheet1 matching the existing table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "test2"
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 44
